# Updated cryptos list — applies the latest price/volume snapshot,
# and shifts the coin ranking list up by one position (LEO drops off
# the bottom of the visible range, Aave is appended as the new #51).
#
# NOTE: Price/Volume columns are stored as *text* in this sheet (some
# prices use "." as a thousands separator, e.g. "28.123.58", which is
# not a valid number) — a leading apostrophe forces Excel to keep the
# literal text instead of silently re-parsing it as a number (which
# would also strip meaningful trailing zeros, e.g. "2.490" -> 2.49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-25: price (D) / volume-1h (E) refresh only ---
$ws.Range("D2").Value = "'28.123.58"
$ws.Range("E2").Value = "'  +0.20%  "
$ws.Range("D3").Value = "'1.873.89"
$ws.Range("E3").Value = "'  -0.02%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("D5").Value = "'313.13"
$ws.Range("E5").Value = "'  +0.06%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  -0.11%  "
$ws.Range("D7").Value = "'0.5139"
$ws.Range("E7").Value = "'  +1.92%  "
$ws.Range("D8").Value = "'0.3891"
$ws.Range("E8").Value = "'  +1.75%  "
$ws.Range("D9").Value = "'0.08395"
$ws.Range("E9").Value = "'  -0.27%  "
$ws.Range("D10").Value = "'1.116"
$ws.Range("D11").Value = "'41.62"
$ws.Range("E11").Value = "'  -0.08%  "
$ws.Range("D12").Value = "'6.206"
$ws.Range("D13").Value = "'20.71"
$ws.Range("E13").Value = "'  +0.97%  "
$ws.Range("D14").Value = "'1.867.10"
$ws.Range("E14").Value = "'  -0.84%  "
$ws.Range("D15").Value = "'7.294"
$ws.Range("E15").Value = "'  +1.33%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "'  -0.19%  "
$ws.Range("E17").Value = "'  +1.09%  "
$ws.Range("D18").Value = "'91.01"
$ws.Range("E18").Value = "'  +0.00%  "
$ws.Range("D19").Value = "'0.06651"
$ws.Range("E19").Value = "'  -0.10%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "'  -1.59%  "
$ws.Range("E21").Value = "'  -0.03%  "
$ws.Range("D22").Value = "'6.051"
$ws.Range("E22").Value = "'  -0.27%  "
$ws.Range("D23").Value = "'28.170.92"
$ws.Range("E23").Value = "'  +0.23%  "
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = "'  -0.26%  "
$ws.Range("D25").Value = "'2.251"
$ws.Range("E25").Value = "'  -0.75%  "

# --- Rows 26-51: ranking list shifts up by one (LEO falls out of
#     range, Aave is newly appended as rank 51) together with a
#     price/volume refresh for every coin that moved up ---
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "'2.083.43"
$ws.Range("E26").Value = "'  -0.77%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.490"
$ws.Range("E27").Value = "'  -3.54%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'158.33"
$ws.Range("E28").Value = "'  +0.84%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.62"
$ws.Range("E29").Value = "'  -0.04%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'125.13"
$ws.Range("E30").Value = "'  -0.74%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1062"
$ws.Range("E31").Value = "'  +1.26%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.040"
$ws.Range("E32").Value = "'  -0.81%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.897"
$ws.Range("E33").Value = "'  +5.03%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.598"
$ws.Range("E34").Value = "'  -0.41%  "

$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.742"
$ws.Range("E35").Value = "'  +0.54%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02439"
$ws.Range("E36").Value = "'  -0.41%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06537"
$ws.Range("E37").Value = "'  -0.06%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2187"
$ws.Range("E38").Value = "'  +0.90%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.207"
$ws.Range("E39").Value = "'  -0.94%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6517"
$ws.Range("E40").Value = "'  +0.05%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "'5.023"
$ws.Range("E41").Value = "'  +2.71%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.229"
$ws.Range("E42").Value = "'  -1.36%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.32"
$ws.Range("E43").Value = "'  -0.10%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6101"
$ws.Range("E44").Value = "'  -1.44%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.04"
$ws.Range("E45").Value = "'  -0.13%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.280"
$ws.Range("E46").Value = "'  -1.72%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.675"
$ws.Range("E47").Value = "'  -0.21%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.011"
$ws.Range("E48").Value = "'  -0.04%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.217"
$ws.Range("E49").Value = "'  -0.14%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'121.65"
$ws.Range("E50").Value = "'  +0.62%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'77.90"
$ws.Range("E51").Value = "'  -2.98%  "
